$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.68%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.84%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.893"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.11%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06423"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.59%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.51%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.242"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-7.30%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8823"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.17%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1513"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.62%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05037"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07504"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02913"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.68%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08995"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.49%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001566"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006408"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.33%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005696"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.98%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.09%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.315"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.26%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.00%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1336"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.61%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.908"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.27%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.42%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.40%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.15%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.02%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "13.96%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04140"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.80%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006819"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.38%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.19%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "13.83%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01170"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.58%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005199"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.42%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.487"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-22.20%"
